$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.447.29"
$ws.Range("E2").Value = "  -3.96%  "

# Row 3
$ws.Range("D3").Value = "'1.948.94"
$ws.Range("E3").Value = "  -2.86%  "

# Row 4
$ws.Range("D4").Value = "'1.017"
$ws.Range("E4").Value = "  +0.50%  "

# Row 5
$ws.Range("D5").Value = "'321.03"
$ws.Range("E5").Value = "  -2.69%  "

# Row 6
$ws.Range("E6").Value = "  +0.33%  "

# Row 7
$ws.Range("D7").Value = "'0.4753"
$ws.Range("E7").Value = "  -5.28%  "

# Row 8
$ws.Range("D8").Value = "'0.4026"
$ws.Range("E8").Value = "  -4.89%  "

# Row 9
$ws.Range("D9").Value = "'53.63"
$ws.Range("E9").Value = "  -0.84%  "

# Row 10
$ws.Range("D10").Value = "'0.08476"
$ws.Range("E10").Value = "  -6.12%  "

# Row 11
$ws.Range("E11").Value = "  -5.57%  "

# Row 12
$ws.Range("D12").Value = "'21.90"
$ws.Range("E12").Value = "  -6.29%  "

# Row 13
$ws.Range("D13").Value = "'1.984.74"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14
$ws.Range("D14").Value = "'7.596"
$ws.Range("E14").Value = "  -6.03%  "

# Row 15
$ws.Range("D15").Value = "'6.191"
$ws.Range("E15").Value = "  -4.73%  "

# Row 16
$ws.Range("D16").Value = "'1.017"
$ws.Range("E16").Value = "  +0.39%  "

# Row 17
$ws.Range("D17").Value = "'0.00001076"
$ws.Range("E17").Value = "  -3.45%  "

# Row 18
$ws.Range("D18").Value = "'88.68"
$ws.Range("E18").Value = "  -6.02%  "

# Row 19
$ws.Range("D19").Value = "'0.06636"
$ws.Range("E19").Value = "  -0.40%  "

# Row 20
$ws.Range("D20").Value = "'18.60"
$ws.Range("E20").Value = "  -5.63%  "

# Row 21
$ws.Range("D21").Value = "'1.013"
$ws.Range("E21").Value = "  +0.18%  "

# Row 22
$ws.Range("D22").Value = "'5.808"
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("D23").Value = "'28.506.55"
$ws.Range("E23").Value = "  -3.81%  "

# Row 24
$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = "  -4.44%  "

# Row 25
$ws.Range("D25").Value = "'2.297"
$ws.Range("E25").Value = "  -0.13%  "

# Row 26
$ws.Range("D26").Value = "'2.167.28"
$ws.Range("E26").Value = "  -3.96%  "

# Row 27
$ws.Range("D27").Value = "'153.88"
$ws.Range("E27").Value = "  -3.11%  "

# Row 28
$ws.Range("E28").Value = "  -3.06%  "

# Row 29
$ws.Range("D29").Value = "'5.917"
$ws.Range("E29").Value = "  -7.55%  "

# Row 30
$ws.Range("D30").Value = "'2.151"
$ws.Range("E30").Value = "  -6.93%  "

# Row 31
$ws.Range("D31").Value = "'123.52"
$ws.Range("E31").Value = "  -3.72%  "

# Row 32
$ws.Range("D32").Value = "'0.9935"
$ws.Range("E32").Value = "  -6.13%  "

# Row 33
$ws.Range("D33").Value = "'0.09559"
$ws.Range("E33").Value = "  -4.15%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.604"
$ws.Range("E34").Value = "  -4.22%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.672"
$ws.Range("E35").Value = "  -3.21%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.434"
$ws.Range("E36").Value = "  -8.91%  "

# Row 37
$ws.Range("D37").Value = "'0.02329"

# Row 38
$ws.Range("D38").Value = "'0.06208"
$ws.Range("E38").Value = "  -2.72%  "

# Row 39
$ws.Range("D39").Value = "'1.256"
$ws.Range("E39").Value = "  -4.20%  "

# Row 40
$ws.Range("D40").Value = "'8.713"
$ws.Range("E40").Value = "  -6.61%  "

# Row 41
$ws.Range("D41").Value = "'0.6219"
$ws.Range("E41").Value = "  -5.35%  "

# Row 42
$ws.Range("D42").Value = "'11.07"
$ws.Range("E42").Value = "  -5.65%  "

# Row 43
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("D44").Value = "'0.1920"
$ws.Range("E44").Value = "  -6.65%  "

# Row 45
$ws.Range("D45").Value = "'1.333"
$ws.Range("E45").Value = "  +1.97%  "

# Row 46
$ws.Range("D46").Value = "'0.5938"
$ws.Range("E46").Value = "  -6.70%  "

# Row 47
$ws.Range("D47").Value = "'12.97"
$ws.Range("E47").Value = "  -3.40%  "

# Row 48
$ws.Range("D48").Value = "'2.059"
$ws.Range("E48").Value = "  -6.31%  "

# Row 49
$ws.Range("D49").Value = "'3.399"
$ws.Range("E49").Value = "  -3.18%  "

# Row 50
$ws.Range("D50").Value = "'0.00000000333"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51
$ws.Range("D51").Value = "'0.06810"
$ws.Range("E51").Value = "  -2.59%  "
